# feat: add 2022-Q1 data
#
# 1) Insert a new worksheet "2022-Q1" right before "总计", carrying the
#    per-fund holding detail (same layout as the other quarter sheets).
# 2) Prepend a new summary row for "2022-Q1" (count=17, value=8.1 yi) at
#    the top of the "总计" sheet's data, shifting the existing rows down.

$wb = $excel.ActiveWorkbook

# NOTE: `Worksheets.Add($before)` re-purposes the handle that is passed
# in as the `$before` argument to BECOME the handle for the freshly
# created sheet (this COM shim does not hand back a separate object for
# "the sheet that used to be at that slot"). So we must not keep a
# variable pointing at "总计" across the Add() call - the sheet is
# looked up again, by a fresh call, once the new sheet is in place.

$sheetCountBeforeAdd = $wb.Worksheets.Count
$newSheet = $wb.Worksheets.Add($wb.Worksheets.Item($sheetCountBeforeAdd))
$newSheet.Name = "2022-Q1"

# A plain quarterly sheet to copy the header/index-column look from.
$templateSheet = $wb.Worksheets.Item(1 + 1)

# Per-fund detail rows: code, name, fund scale, total stock position,
# position ratio, holding value (yi), position rank.
$fundData = @(
    @("161834","银华鑫锐灵活配置混合（LOF）","67.33","81.90","2.56","1.7236",4),
    @("501022","银华鑫盛灵活配置混合（LOF）","61.98","79.75","2.44","1.5123",4),
    @("550015","信诚至远灵活配置混合A","41.04","79.97","2.69","1.1040",6),
    @("014677","中信保诚至远动力混合E","41.04","79.97","2.69","1.1040",6),
    @("001736","圆信永丰优加生活股票","54.81","82.64","1.56","0.8550",9),
    @("240008","华宝收益增长混合","9.38","93.55","3.68","0.3452",9),
    @("012370","银华鑫利一年持有期混合型证券投资基金","11.01","80.06","2.27","0.2499",6),
    @("009913","中信保诚成长动力混合","8.86","80.14","2.70","0.2392",6),
    @("004959","圆信永丰优悦生活混合","13.70","79.90","1.54","0.2110",10),
    @("008245","圆信永丰致优混合A","12.80","82.48","1.56","0.1997",10),
    @("550016","信诚至远灵活配置混合C","6.06","79.97","2.69","0.1630",6),
    @("001581","华安沪港深通精选灵活配置混合","4.92","92.91","3.12","0.1535",10),
    @("008246","圆信永丰致优混合C","5.11","82.48","1.56","0.0797",10),
    @("560660","新华中证云计算50交易型开放式指数证券投资基金","1.75","97.21","3.06","0.0536",8),
    @("001707","诺安高端制造股票","1.35","91.74","3.50","0.0472",7),
    @("001534","华宝万物互联灵活配置混合","1.06","92.81","3.91","0.0414",5),
    @("002152","华宝核心优势灵活配置混合","0.45","90.91","3.76","0.0169",4)
)

$lastRow = 1 + $fundData.Length

# Header row (B1:H1) - copy the bold/centered/bordered look from an
# existing quarter sheet's header.
$templateSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$newSheet.Cells.Item(1,2).Value2 = "基金代码"
$newSheet.Cells.Item(1,3).Value2 = "基金名称"
$newSheet.Cells.Item(1,4).Value2 = "基金规模"
$newSheet.Cells.Item(1,5).Value2 = "股票总仓位"
$newSheet.Cells.Item(1,6).Value2 = "仓位占比"
$newSheet.Cells.Item(1,7).Value2 = "持有市值(亿元)"
$newSheet.Cells.Item(1,8).Value2 = "仓位排名"

# Column A index cells (A2:A{lastRow}) - copy the bold/centered/bordered
# look used by every other quarter sheet's index column.
$templateSheet.Range("A2").Copy()
$newSheet.Range(("A2:A" + $lastRow)).PasteSpecial(-4122)

# Force columns B:G to be read back as text (fund codes such as
# "001736" must keep their leading zeros, and the numeric-looking
# scale/position/ratio/value columns are stored as text in every other
# quarter sheet too) - then drop the temporary format so the cells end
# up with the plain, unstyled look the other sheets use.
$dataRange = $newSheet.Range(("B2:G" + $lastRow))
$dataRange.NumberFormat = "@"

$r = 2
foreach ($row in $fundData) {
    $newSheet.Cells.Item($r,1).Value2 = $r - 2
    $newSheet.Cells.Item($r,2).Value2 = $row[0]
    $newSheet.Cells.Item($r,3).Value2 = $row[1]
    $newSheet.Cells.Item($r,4).Value2 = $row[2]
    $newSheet.Cells.Item($r,5).Value2 = $row[3]
    $newSheet.Cells.Item($r,6).Value2 = $row[4]
    $newSheet.Cells.Item($r,7).Value2 = $row[5]
    $newSheet.Cells.Item($r,8).Value2 = $row[6]
    $r = $r + 1
}

$dataRange.ClearFormats()

# Match the page margins used by the other quarter sheets (0.75in /
# 1in / 0.5in).
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

$newSheet.Range("A1").Select()

# ---------------------------------------------------------------------
# 2) Prepend the "2022-Q1" summary row to the "总计" sheet (now the
#    last sheet again), shifting the existing rows down by one.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

$oldB = @()
$oldC = @()
$oldD = @()
for ($row = 2; $row -le 6; $row++) {
    $oldB += $totalSheet.Cells.Item($row, 2).Value2
    $oldC += $totalSheet.Cells.Item($row, 3).Value2
    $oldD += $totalSheet.Cells.Item($row, 4).Value2
}

# Make sure the newly-exposed row (7) has the same style as the other
# index cells in column A before writing into it.
$totalSheet.Range("A6").Copy()
$totalSheet.Range("A7").PasteSpecial(-4122)

for ($i = $oldB.Length - 1; $i -ge 0; $i--) {
    $destRow = $i + 3
    $totalSheet.Cells.Item($destRow, 1).Value2 = $i + 1
    $totalSheet.Cells.Item($destRow, 2).Value2 = $oldB[$i]
    $totalSheet.Cells.Item($destRow, 3).Value2 = $oldC[$i]
    $totalSheet.Cells.Item($destRow, 4).Value2 = $oldD[$i]
}

$totalSheet.Cells.Item(2, 1).Value2 = 0
$totalSheet.Cells.Item(2, 2).Value2 = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value2 = 17
$totalSheet.Cells.Item(2, 4).Value2 = 8.1
